$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 7 (only header row 1 and data row 2 remain)
$ws.Range("A3:G7").EntireRow.Delete()

# Update row 2 values
$ws.Range("A2").Value = "Daniel"
$ws.Range("B2").Value = "empty"
$ws.Range("C2").Value = "empty"
$ws.Range("D2").Value = 0
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "empty"
$ws.Range("G2").Value = 0
